$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text fixes (Volume/Number, report week dates) ---
$volCell = $ws.Cells.Item(8, 1)
$volCell.Characters(21, 2).Text = "19"

$dateCell = $ws.Cells.Item(9, 3)
$dateCell.Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 4
$ws.Range("K15").Value = -20
$ws.Range("L15").Value = -60
$ws.Range("M15").Value = 33.333333333333
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 58
$ws.Range("J16").Value = 77
$ws.Range("K16").Value = -24.675324675324
$ws.Range("L16").Value = 81.25
$ws.Range("M16").Value = -36.263736263736
$ws.Range("N16").Value = -86.320754716981
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -77.777777777777
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 28.571428571428
$ws.Range("I17").Value = 136
$ws.Range("J17").Value = 118
$ws.Range("K17").Value = 15.254237288135
$ws.Range("L17").Value = 65.853658536585
$ws.Range("M17").Value = 81.333333333333
$ws.Range("N17").Value = 12.396694214876
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 250
$ws.Range("F18").Value = 20
$ws.Range("H18").Value = 122.222222222222
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 44
$ws.Range("K18").Value = 27.272727272727
$ws.Range("L18").Value = 36.585365853658
$ws.Range("M18").Value = -41.666666666666
$ws.Range("N18").Value = -88.501026694045
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = -44.444444444444
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -46.428571428571
$ws.Range("I19").Value = 108
$ws.Range("J19").Value = 144
$ws.Range("K19").Value = -25
$ws.Range("L19").Value = 30.12048192771
$ws.Range("M19").Value = -7.692307692307
$ws.Range("N19").Value = -43.455497382199
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 8
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -10
$ws.Range("I20").Value = 94
$ws.Range("J20").Value = 89
$ws.Range("K20").Value = 5.617977528089
$ws.Range("L20").Value = 46.875
$ws.Range("M20").Value = 6.818181818181
$ws.Range("N20").Value = -92.857142857142
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -23.684210526315
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = 2.061855670103
$ws.Range("I21").Value = 458
$ws.Range("J21").Value = 478
$ws.Range("K21").Value = -4.18410041841
$ws.Range("L21").Value = 45.396825396825
$ws.Range("M21").Value = -2.966101694915
$ws.Range("N21").Value = -82.074363992172
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -55.555555555555
$ws.Range("L22").Value = 33.333333333333
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -55.263157894736
$ws.Range("F24").Value = 113
$ws.Range("G24").Value = 125
$ws.Range("H24").Value = -9.6
$ws.Range("I24").Value = 536
$ws.Range("J24").Value = 519
$ws.Range("K24").Value = 3.275529865125
$ws.Range("L24").Value = 43.31550802139
$ws.Range("M24").Value = 89.399293286219
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 32
$ws.Range("H25").Value = -3.030303030303
$ws.Range("I25").Value = 194
$ws.Range("J25").Value = 174
$ws.Range("K25").Value = 11.494252873563
$ws.Range("L25").Value = 43.703703703703
$ws.Range("M25").Value = -7.619047619047
$ws.Range("C26").Value = 3
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 100
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -14.285714285714
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 10.526315789473
$ws.Range("L27").Value = 16.666666666666

# --- Row 22: text placeholders -> real numbers (match numeric style by format) ---
$d22 = $ws.Range("D22"); $d22.NumberFormat = "#,##0"; $d22.Value = 2
$e22 = $ws.Range("E22"); $e22.NumberFormat = '#,##0.0;"-"#,##0.0'; $e22.Value = -100
$g22 = $ws.Range("G22"); $g22.NumberFormat = "#,##0"; $g22.Value = 2
$h22 = $ws.Range("H22"); $h22.NumberFormat = '#,##0.0;"-"#,##0.0'; $h22.Value = -100

# --- Row 27: text placeholder "0" -> real number 2 ---
$c27 = $ws.Range("C27"); $c27.NumberFormat = "#,##0"; $c27.Value = 2

# --- Rows 28-30: numbers -> text placeholders ("0" / "***.*") ---
$g28 = $ws.Range("G28"); $g28.NumberFormat = "@"; $g28.Value = "0"
$h28 = $ws.Range("H28"); $h28.NumberFormat = "@"; $h28.Value = "***.*"
$g29 = $ws.Range("G29"); $g29.NumberFormat = "@"; $g29.Value = "0"
$h29 = $ws.Range("H29"); $h29.NumberFormat = "@"; $h29.Value = "***.*"
$f30 = $ws.Range("F30"); $f30.NumberFormat = "@"; $f30.Value = "0"
$g30 = $ws.Range("G30"); $g30.NumberFormat = "@"; $g30.Value = "0"
$h30 = $ws.Range("H30"); $h30.NumberFormat = "@"; $h30.Value = "***.*"
